$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 previously held a placeholder number (0) styled with a bold/bordered/
# centered xf. The new content is the full questions payload (now pretty-
# printed JSON instead of a Python-dict-literal string), stored as a plain
# un-styled cell. Clear the old style first so the long pasted value does not
# inherit the prior font/border/alignment.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = "questions = [`n    {`n        `"title`": `"You want to display a page that has the path \u201chttps://domain_name/company/history/last_year.\u201d Which of the following file structures do you need to create in order to display that page?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"src/pages/company/history/last_year.js`",`n            `"src/company/history/last_year.js`",`n            `"src/templates/company/history/last_year.js`",`n            `"src/pages/last_year.js`"`n        ],`n        `"score`": `"src/pages/company/history/last_year.js`"`n    },`n    {`n        `"title`": `"You need to display a product\u2019s dynamic image from a content management system (CMS). Which of the following pieces of code will allow you to display the image?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"&ltGatsbyImage image={data.product.avatar} alt={data.product.name} /&gt`",`n            `"&ltGatsbyImage image={data.product.avatar.childImageSharp.gatsbyImageData} alt={data.product.name} /&gt`",`n            `"&ltStaticImage image={data.product.avatar} alt={data.product.name} /&gt`",`n            `"&ltStaticImage image={data.product.avatar.childImageSharp.gatsbyImageData} alt={data.product.name} /&gt`"`n        ],`n        `"score`": `"&ltGatsbyImage image={data.product.avatar.childImageSharp.gatsbyImageData} alt={data.product.name} /&gt`"`n    },`n    {`n        `"title`": `"You use a content management system (CMS) to obtain blog posts and need to create a page for each post. You have written the code shown below to create the pages. Which of the following path structures will each page have?  result.data.allPosts.edges.forEach(({ node }) =&gt {\n     createPage({\n       path: ``/`${node.url}``,\n       component: path.resolve(``src/templates/post.js``),\n       context: {\n         url: node.url,\n       },\n     })\n   })`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"https://domain_name/allPosts/post_url`",`n            `"https://domain_name/posts/post_url`",`n            `"https://domain_name/templates/post_url`",`n            `"https://domain_name /post_url`"`n        ],`n        `"score`": `"https://domain_name /post_url`"`n    },`n    {`n        `"title`": `"True or false: \u201cgatsby new [&ltsite-name&gt [&ltstarter-url&gt]]\u201d can create a site from a starter.`",`n        `"ques_type`": 11,`n        `"options`": [`n            `"true`",`n            `"false`"`n        ],`n        `"score`": `"True`"`n    }`n]"

# Restore the row to its natural (non-custom) height after the multi-line
# text assignment would otherwise pin an explicit autofit height.
$ws.Rows(1).AutoFit()

# A2 (the duplicate/old shared-string cell) is removed entirely so the sheet
# only has the single A1 cell left.
$ws.Range("A2").ClearContents()
